{"js": "// Insert a new bullet paragraph \"-T\u1ea3i d\u1eef li\u1ec7u product l\u00ean admin\" right\n// after the existing \"-T\u1ea3i d\u1eef li\u1ec7u c\u1ee7a Category l\u00ean trang admin.\" paragraph\n// in the \"B\u00e1o c\u00e1o tu\u1ea7n 6\" section.\n\nconst body = context.document.body;\nconst results = body.search(\"-T\u1ea3i d\u1eef li\u1ec7u c\u1ee7a Category l\u00ean trang admin.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\nanchorParagraph.load(\"text\");\nawait context.sync();\n\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"-T\u1ea3i d\u1eef li\u1ec7u product l\u00ean admin\",\n  \"After\"\n);\n\n// Match the formatting used by the surrounding report bullets: Times New\n// Roman font on the run (the paragraph mark itself also carries this font\n// in the source XML, which insertParagraph's clone of the anchor paragraph\n// already provides).\nnewParagraph.font.set({ name: \"Times New Roman\" });\n\nawait context.sync();\n", "ps1": "# Insert a new bullet paragraph \"-T\u1ea3i d\u1eef li\u1ec7u product l\u00ean admin\" right\n# after the existing \"-T\u1ea3i d\u1eef li\u1ec7u c\u1ee7a Category l\u00ean trang admin.\" paragraph\n# in the \"B\u00e1o c\u00e1o tu\u1ea7n 6\" weekly report section.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"-T\u1ea3i d\u1eef li\u1ec7u c\u1ee7a Category l\u00ean trang admin.\"\n$newText = \"-T\u1ea3i d\u1eef li\u1ec7u product l\u00ean admin\"\n\n$paras = $d.Paragraphs\n$targetIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq $anchorText) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Anchor paragraph not found\"\n}\n\n$target = $paras.Item($targetIndex)\n$target.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = $newText\n"}
